# NIT-9009260332.xlsx update
# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The list of late-payment periods (Periodo Mora, col E) together with their
# Valor Mora (col F) is reordered (it is now listed from the most recent
# period down to the oldest one, i.e. the previous order is reversed), and
# the Salario Basico (col G) for every one of those rows is refreshed to a
# new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New (period -> valor mora) pairs for rows 16..35, already in the order
# they must appear after the edit (row 16 = newest period, row 35 = oldest).
$periodos = @(
    @{ Row = 16; Periodo = "2009"; ValorMora = 30916 },
    @{ Row = 17; Periodo = "2008"; ValorMora = 33125 },
    @{ Row = 18; Periodo = "2007"; ValorMora = 33125 },
    @{ Row = 19; Periodo = "2006"; ValorMora = 33125 },
    @{ Row = 20; Periodo = "2005"; ValorMora = 33125 },
    @{ Row = 21; Periodo = "2004"; ValorMora = 33125 },
    @{ Row = 22; Periodo = "2003"; ValorMora = 33125 },
    @{ Row = 23; Periodo = "2002"; ValorMora = 33125 },
    @{ Row = 24; Periodo = "2001"; ValorMora = 33125 },
    @{ Row = 25; Periodo = "1912"; ValorMora = 33125 },
    @{ Row = 26; Periodo = "1911"; ValorMora = 33125 },
    @{ Row = 27; Periodo = "1910"; ValorMora = 33125 },
    @{ Row = 28; Periodo = "1909"; ValorMora = 33125 },
    @{ Row = 29; Periodo = "1903"; ValorMora = 26041 },
    @{ Row = 30; Periodo = "1902"; ValorMora = 31249 },
    @{ Row = 31; Periodo = "1901"; ValorMora = 31249 },
    @{ Row = 32; Periodo = "1812"; ValorMora = 31249 },
    @{ Row = 33; Periodo = "1811"; ValorMora = 31249 },
    @{ Row = 34; Periodo = "1810"; ValorMora = 31249 },
    @{ Row = 35; Periodo = "1809"; ValorMora = 31249 }
)

foreach ($p in $periodos) {
    $r = $p.Row
    $ws.Cells.Item($r, 5).Value = $p.Periodo      # column E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $p.ValorMora    # column F - Valor Mora
}

# Salario Basico (column G) is refreshed to the new value for every worker row.
$ws.Range("G16:G35").Value = 828116
